$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the old "code" column (B), pushing
# product/inventary/cost/benefit/price/date to the right.
$ws.Range("B1:E1").EntireColumn.Insert()

# The insert operation made the new data cells (B2:E32) inherit column A's
# style (border/bold) because column A had a style on every data row. The
# new columns should be plain/unstyled like the other data columns, so
# clear that inherited formatting.
$ws.Range("B2:E32").ClearFormats()

# The header row's original A1 cell was blank, so the newly inserted
# header cells did not inherit the bold/bordered header style. Copy the
# formatting from the neighboring header cell (now F1, "code") onto them.
$ws.Range("F1").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new header labels for the inserted columns.
$ws.Range("B1").Value = "Unnamed: 0.3"
$ws.Range("C1").Value = "Unnamed: 0.2"
$ws.Range("D1").Value = "Unnamed: 0.1"
$ws.Range("E1").Value = "Unnamed: 0"

# Populate the new columns with the same running index as column A
# for each data row (rows 2..32).
for ($r = 2; $r -le 32; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $idx
    $ws.Cells.Item($r, 3).Value = $idx
    $ws.Cells.Item($r, 4).Value = $idx
    $ws.Cells.Item($r, 5).Value = $idx
}

# Row 32's "inventary" value (now column H) changes from 10 to 0.
$ws.Range("H32").Value = 0
